$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report was inserted at the top of the "Camote" price
# block (rows 249-250), which pushes all the subsequent rows (old 249-354)
# down by two (to 251-356). Excel's native row-insert handles that shift
# (values, styles, shared-string refs, dimension, etc.) for us.
$ws.Rows("249:250").Insert()

# Fill in the two brand-new rows with this week's "Camote" prices.
$ws.Range("A249").Value = 11
$ws.Range("B249").Value = 'Vega Monumental Concepción'
$ws.Range("C249").Value = 'Bíobío'
$ws.Range("D249").Value = 44992
$ws.Range("E249").Value = 8
$ws.Range("F249").Value = 100112045
$ws.Range("G249").Value = 'Zapallo'
$ws.Range("H249").Value = 'Camote'
$ws.Range("I249").Value = '1a (cosecha)'
$ws.Range("J249").Value = 600
$ws.Range("K249").Value = 600
$ws.Range("L249").Value = 650
$ws.Range("M249").Value = 625
$ws.Range("N249").Value = '$/kilo (volumen en unidades)'
$ws.Range("O249").Value = 'Región Metropolitana'
$ws.Range("P249").Value = 625
$ws.Range("Q249").Value = 1
$ws.Range("R249").Value = 'Hortaliza'

$ws.Range("A250").Value = 11
$ws.Range("B250").Value = 'Vega Monumental Concepción'
$ws.Range("C250").Value = 'Bíobío'
$ws.Range("D250").Value = 44992
$ws.Range("E250").Value = 8
$ws.Range("F250").Value = 100112045
$ws.Range("G250").Value = 'Zapallo'
$ws.Range("H250").Value = 'Camote'
$ws.Range("I250").Value = '2a (cosecha)'
$ws.Range("J250").Value = 300
$ws.Range("K250").Value = 500
$ws.Range("L250").Value = 500
$ws.Range("M250").Value = 500
$ws.Range("N250").Value = '$/kilo (volumen en unidades)'
$ws.Range("O250").Value = 'Región Metropolitana'
$ws.Range("P250").Value = 500
$ws.Range("Q250").Value = 1
$ws.Range("R250").Value = 'Hortaliza'
